$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell C1: "Remark" -> "Expected " (trailing space)
$ws.Range("C1").Value = "Expected "

# Apply bold font + yellow fill to the header row (A1:C1).
# Build the format on A1 first (creates one clean combined style),
# then copy that formatting onto B1 and C1 via copy/paste-special so
# only a single new font/fill/cellXf trio is introduced.
$a1 = $ws.Range("A1")
$a1.Font.Bold = $true
$a1.Interior.Color = 65535

$a1.Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)

# Move/keep the active selection at E8, matching the saved view state.
[void]$ws.Range("E8").Select()

# Print setup: force portrait orientation (adds <pageSetup .../>).
$ws.PageSetup.Orientation = 1

Write-Host "done"
